# Update the date header and the 25 division problems in the practice table.
$d = $word.ActiveDocument

# --- Date heading ---------------------------------------------------------
[void]$d.Content.Find.Execute("2025-07-02 Wednesday", $true, $false, $false,
                               $false, $false, $true, 1, $false,
                               "2025-07-03 Thursday", 2)

# --- Division problems (addressed by table cell so the two cells that ----
# --- originally shared the text "66÷9=" still get their own new value) ---
$t = $d.Tables.Item(1)

$values = @(
    @(1, 1, "93÷2="),
    @(1, 2, "44÷3="),
    @(1, 3, "74÷8="),
    @(1, 4, "30÷4="),
    @(1, 5, "33÷4="),

    @(5, 1, "37÷3="),
    @(5, 2, "46÷8="),
    @(5, 3, "66÷7="),
    @(5, 4, "36÷8="),
    @(5, 5, "38÷7="),

    @(9, 1, "37÷7="),
    @(9, 2, "11÷5="),
    @(9, 3, "33÷6="),
    @(9, 4, "36÷6="),
    @(9, 5, "91÷5="),

    @(13, 1, "89÷7="),
    @(13, 2, "92÷4="),
    @(13, 3, "86÷2="),
    @(13, 4, "93÷2="),
    @(13, 5, "53÷7="),

    @(17, 1, "67÷8="),
    @(17, 2, "66÷2="),
    @(17, 3, "70÷3="),
    @(17, 4, "86÷3="),
    @(17, 5, "36÷2=")
)

foreach ($entry in $values) {
    $row = $entry[0]
    $col = $entry[1]
    $text = $entry[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $text
}
